# Auto-generated script applying scheduled market-data refresh values
# to the Lich_Profits crafting-profit workbook (one worksheet per crafting class).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3465.0715
$ws.Range("I88").Value = 3484
$ws.Range("J88").Value = 3457.5
$ws.Range("K88").Value = 3484
$ws.Range("L88").Value = 3457.5
$ws.Range("M88").Value = -3078
$ws.Range("N88").Value = -4269.5
$ws.Range("H91").Value = 3465.0715
$ws.Range("I91").Value = 3484
$ws.Range("J91").Value = 3457.5
$ws.Range("K91").Value = 3484
$ws.Range("L91").Value = 3457.5
$ws.Range("M91").Value = -2080
$ws.Range("N91").Value = -6265.5
$ws.Range("H132").Value = 1267.5532
$ws.Range("I132").Value = 1014.89746
$ws.Range("J132").Value = 2499.25
$ws.Range("K132").Value = 3044.69238
$ws.Range("L132").Value = 7497.75
$ws.Range("M132").Value = -514.69238
$ws.Range("N132").Value = -12557.75
$ws.Range("H136").Value = 182442.25
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 182442.25
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 182442.25
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -192642.25
$ws.Range("H137").Value = 56361.81
$ws.Range("I137").Value = 67105.586
$ws.Range("K137").Value = 201316.758
$ws.Range("M137").Value = -198766.758

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16119
$ws.Range("I2").Value = 18590
$ws.Range("K2").Value = 18590
$ws.Range("M2").Value = -18477
$ws.Range("H61").Value = 4593.4517
$ws.Range("I61").Value = 2946.4783
$ws.Range("K61").Value = 2946.4783
$ws.Range("M61").Value = -2734.4783
$ws.Range("H74").Value = 43346.812
$ws.Range("I74").Value = 43346.812
$ws.Range("K74").Value = 43346.812
$ws.Range("M74").Value = -42472.812
$ws.Range("H77").Value = 43346.812
$ws.Range("I77").Value = 43346.812
$ws.Range("K77").Value = 216734.06
$ws.Range("M77").Value = -212366.06
$ws.Range("H88").Value = 923.17645
$ws.Range("I88").Value = 943.9286
$ws.Range("J88").Value = 908.65
$ws.Range("K88").Value = 943.9286
$ws.Range("L88").Value = 908.65
$ws.Range("M88").Value = -537.9286
$ws.Range("N88").Value = -1720.65
$ws.Range("H91").Value = 923.17645
$ws.Range("I91").Value = 943.9286
$ws.Range("J91").Value = 908.65
$ws.Range("K91").Value = 943.9286
$ws.Range("L91").Value = 908.65
$ws.Range("M91").Value = 460.0714
$ws.Range("N91").Value = -3716.65
$ws.Range("H97").Value = 1593.3793
$ws.Range("I97").Value = 1108.8636
$ws.Range("K97").Value = 1108.8636
$ws.Range("M97").Value = -612.8635999999999
$ws.Range("H112").Value = 40257.668
$ws.Range("J112").Value = 40257.668
$ws.Range("L112").Value = 40257.668
$ws.Range("N112").Value = -43211.668
$ws.Range("H113").Value = 79000
$ws.Range("J113").Value = 79000
$ws.Range("L113").Value = 79000
$ws.Range("N113").Value = -87678
$ws.Range("H116").Value = 16119
$ws.Range("I116").Value = 18590
$ws.Range("K116").Value = 18590
$ws.Range("M116").Value = -16296
$ws.Range("H122").Value = 6501.136
$ws.Range("I122").Value = 4211.2354
$ws.Range("K122").Value = 12633.7062
$ws.Range("M122").Value = -10183.7062
$ws.Range("H136").Value = 4593.4517
$ws.Range("I136").Value = 2946.4783
$ws.Range("K136").Value = 8839.4349
$ws.Range("M136").Value = -6289.4349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16119
$ws.Range("I3").Value = 18590
$ws.Range("K3").Value = 18590
$ws.Range("M3").Value = -18476
$ws.Range("H86").Value = 1829.2941
$ws.Range("I86").Value = 1829.2941
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1829.2941
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -706.2941000000001
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1829.2941
$ws.Range("I89").Value = 1829.2941
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9146.470499999999
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -3530.470499999999
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 8263.723
$ws.Range("I99").Value = 8632.177
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 8632.177
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -7134.177
$ws.Range("N99").Value = -4996
$ws.Range("H141").Value = 53332.89
$ws.Range("J141").Value = 53332.89
$ws.Range("L141").Value = 53332.89
$ws.Range("N141").Value = -63692.89

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 432776.56
$ws.Range("I31").Value = 590589.6
$ws.Range("J31").Value = 97423.75
$ws.Range("K31").Value = 590589.6
$ws.Range("L31").Value = 97423.75
$ws.Range("M31").Value = -590294.6
$ws.Range("N31").Value = -98013.75
$ws.Range("H34").Value = 432776.56
$ws.Range("I34").Value = 590589.6
$ws.Range("J34").Value = 97423.75
$ws.Range("K34").Value = 590589.6
$ws.Range("L34").Value = 97423.75
$ws.Range("M34").Value = -590387.6
$ws.Range("N34").Value = -97827.75
$ws.Range("H58").Value = 2026.1818
$ws.Range("I58").Value = 2070.05
$ws.Range("K58").Value = 2070.05
$ws.Range("M58").Value = -1867.05
$ws.Range("H99").Value = 772306.3
$ws.Range("I99").Value = 1002908.4
$ws.Range("K99").Value = 1002908.4
$ws.Range("M99").Value = -1001410.4
$ws.Range("H126").Value = 772306.3
$ws.Range("I126").Value = 1002908.4
$ws.Range("K126").Value = 3008725.2
$ws.Range("M126").Value = -3006255.2
$ws.Range("H134").Value = 4336.769
$ws.Range("I134").Value = 4796.613
$ws.Range("J134").Value = 2554.875
$ws.Range("K134").Value = 14389.839
$ws.Range("L134").Value = 7664.625
$ws.Range("M134").Value = -11854.839
$ws.Range("N134").Value = -12734.625
$ws.Range("H136").Value = 2026.1818
$ws.Range("I136").Value = 2070.05
$ws.Range("K136").Value = 6210.150000000001
$ws.Range("M136").Value = -3660.150000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1882.8096
$ws.Range("I39").Value = 682.3333
$ws.Range("J39").Value = 2082.889
$ws.Range("K39").Value = 2046.9999
$ws.Range("L39").Value = 6248.667
$ws.Range("M39").Value = -1752.9999
$ws.Range("N39").Value = -6836.667
$ws.Range("H107").Value = 490.05554
$ws.Range("J107").Value = 479.26086
$ws.Range("L107").Value = 1437.78258
$ws.Range("N107").Value = -5277.78258
$ws.Range("H121").Value = 2787.9285
$ws.Range("I121").Value = 1438.625
$ws.Range("K121").Value = 4315.875
$ws.Range("M121").Value = -3005.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7560.4287
$ws.Range("I126").Value = 7226.125
$ws.Range("K126").Value = 21678.375
$ws.Range("M126").Value = -19208.375
$ws.Range("H132").Value = 59518.45
$ws.Range("I132").Value = 82160.16
$ws.Range("J132").Value = 17469.572
$ws.Range("K132").Value = 246480.48
$ws.Range("L132").Value = 52408.716
$ws.Range("M132").Value = -243950.48
$ws.Range("N132").Value = -57468.716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1456.8182
$ws.Range("I22").Value = 1554.1666
$ws.Range("J22").Value = 1340
$ws.Range("K22").Value = 1554.1666
$ws.Range("L22").Value = 1340
$ws.Range("M22").Value = -1259.1666
$ws.Range("N22").Value = -1930
$ws.Range("H27").Value = 1456.8182
$ws.Range("I27").Value = 1554.1666
$ws.Range("J27").Value = 1340
$ws.Range("K27").Value = 1554.1666
$ws.Range("L27").Value = 1340
$ws.Range("M27").Value = -1447.1666
$ws.Range("N27").Value = -1554
$ws.Range("H93").Value = 809.9167
$ws.Range("I93").Value = 824.6818
$ws.Range("J93").Value = 647.5
$ws.Range("K93").Value = 824.6818
$ws.Range("L93").Value = 647.5
$ws.Range("M93").Value = 423.3182
$ws.Range("N93").Value = -3143.5
$ws.Range("H110").Value = 49571.75
$ws.Range("J110").Value = 49571.75
$ws.Range("L110").Value = 49571.75
$ws.Range("N110").Value = -57751.75
$ws.Range("H136").Value = 3280.3928
$ws.Range("I136").Value = 3282.7307
$ws.Range("K136").Value = 9848.1921
$ws.Range("M136").Value = -7298.1921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 705.75
$ws.Range("I107").Value = 660.4545000000001
$ws.Range("K107").Value = 1981.3635
$ws.Range("M107").Value = -61.36350000000016
$ws.Range("H114").Value = 48000
$ws.Range("I114").Value = 48000
$ws.Range("K114").Value = 48000
$ws.Range("M114").Value = -43661
$ws.Range("H122").Value = 3854.3438
$ws.Range("I122").Value = 3782.5386
$ws.Range("J122").Value = 4165.5
$ws.Range("K122").Value = 11347.6158
$ws.Range("L122").Value = 12496.5
$ws.Range("M122").Value = -8897.6158
$ws.Range("N122").Value = -17396.5
$ws.Range("H136").Value = 314036.62
$ws.Range("I136").Value = 358727.6
$ws.Range("K136").Value = 1076182.8
$ws.Range("M136").Value = -1073632.8
